$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''65.401.16'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.21%  '

# Row 3
$ws.Range("D3").Value = '''3.546.04'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +4.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '''600.25'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.80%  '

# Row 6
$ws.Range("D6").Value = '''138.34'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.96%  '

# Row 7
$ws.Range("D7").Value = '''3.545.94'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("E9").Value = '  -0.40%  '

# Row 10
$ws.Range("E10").Value = '  +3.77%  '

# Row 11
$ws.Range("D11").Value = '''6.90'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.96%  '

# Row 12
$ws.Range("D12").Value = '''0.386'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.54%  '

# Row 13
$ws.Range("D13").Value = '''4.150.67'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.11%  '

# Row 14
$ws.Range("E14").Value = '  +3.36%  '

# Row 15
$ws.Range("D15").Value = '''27.21'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +5.02%  '

# Row 16
$ws.Range("D16").Value = '''3.556.20'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.21%  '

# Row 17
$ws.Range("E17").Value = '  +1.43%  '

# Row 18
$ws.Range("D18").Value = '''65.316.03'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.03%  '

# Row 19
$ws.Range("D19").Value = '''10.22'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.60%  '

# Row 20
$ws.Range("D20").Value = '''5.94'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.75%  '

# Row 21
$ws.Range("D21").Value = '''14.29'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.74%  '

# Row 22
$ws.Range("D22").Value = '''393.32'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.16%  '

# Row 23
$ws.Range("D23").Value = '''0.574'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.96%  '

# Row 24
$ws.Range("D24").Value = '''3.693.63'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.10%  '

# Row 25
$ws.Range("D25").Value = '''73.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.51%  '

# Row 26
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("E27").Value = '  +10.90%  '

# Row 28
$ws.Range("D28").Value = '''7.81'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +10.86%  '

# Row 29
$ws.Range("E29").Value = '  -0.28%  '

# Row 30
$ws.Range("E30").Value = '  +3.67%  '

# Row 31
$ws.Range("E31").Value = '  +1.05%  '

# Row 32
$ws.Range("D32").Value = '''3.560.18'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +4.28%  '

# Row 34
$ws.Range("D34").Value = '''23.87'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +4.82%  '

# Row 35
$ws.Range("D35").Value = '''0.145'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.91%  '

# Row 36
$ws.Range("E36").Value = '  +14.56%  '

# Row 37
$ws.Range("D37").Value = '''6.96'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.16%  '

# Row 38
$ws.Range("D38").Value = '''169.55'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.93%  '

# Row 39
$ws.Range("E39").Value = '  +8.07%  '

# Row 40
$ws.Range("D40").Value = '''5.00'
$ws.Range("D40").ClearFormats()

# Row 41
$ws.Range("D41").Value = '''0.0803'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.24%  '

# Row 42
$ws.Range("D42").Value = '''0.827'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.77%  '

# Row 43
$ws.Range("D43").Value = '''26.73'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +20.61%  '

# Row 44
$ws.Range("D44").Value = '''42.62'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.05%  '

# Row 45
$ws.Range("E45").Value = '  -0.07%  '

# Row 46
$ws.Range("D46").Value = '''4.43'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.53%  '

# Row 47
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '''1.68'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +5.66%  '

# Row 48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '''1.20'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +10.45%  '

# Row 49
$ws.Range("E49").Value = '  +4.86%  '

# Row 50
$ws.Range("D50").Value = '''2.405.43'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +10.53%  '

# Row 51
$ws.Range("D51").Value = '''307.05'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +9.89%  '
